# Fix CSEI93 source codes
# Update individual question answers (column B) on the "question_answers" sheet,
# and the corresponding computed totals (column B) on the "outputs" sheet.

$wb = $excel.ActiveWorkbook

$wsQA = $wb.Worksheets.Item("question_answers")
$wsOut = $wb.Worksheets.Item("outputs")

# question_answers sheet: cell -> new value (inline string "1" or "2")
$qaUpdates = @{
    "B3"  = "2"
    "B6"  = "2"
    "B7"  = "2"
    "B11" = "1"
    "B14" = "2"
    "B16" = "2"
    "B18" = "1"
    "B20" = "1"
    "B21" = "1"
    "B22" = "1"
    "B27" = "1"
    "B28" = "2"
    "B31" = "1"
    "B32" = "2"
    "B33" = "1"
    "B36" = "2"
    "B37" = "2"
    "B40" = "1"
    "B41" = "2"
    "B42" = "2"
    "B46" = "2"
    "B49" = "2"
    "B50" = "1"
    "B52" = "1"
    "B55" = "2"
    "B56" = "2"
    "B57" = "2"
    "B59" = "1"
}

foreach ($cell in $qaUpdates.Keys) {
    $rng = $wsQA.Range($cell)
    # Force text storage so the numeric-looking answer code ("1"/"2") is
    # kept as a string, matching the existing cells in this column.
    $rng.NumberFormat = "@"
    $rng.Value = $qaUpdates[$cell]
}

# outputs sheet: cell -> new numeric value
$outUpdates = @{
    "B2" = 18
    "B3" = 5
    "B4" = 4
    "B5" = 5
    "B7" = 32
}

foreach ($cell in $outUpdates.Keys) {
    $wsOut.Range($cell).Value = $outUpdates[$cell]
}
